# "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The "Periodo Mora" column (E16:E23) lists the overdue periods for each
# worker. The oldest periods are dropped and new ones are appended, which
# (together with a couple of matching "Valor Mora" amounts in column F)
# shifts every period label one slot newer:
#
#   old: 1803, 2308, 2309, 2310, 2311, 2312, 2401, 2402
#   new: 1803, 2402, 2401, 2312, 2311, 2310, 2309, 2308
#
# and the "Valor Mora" values that travelled with the 2308/2402 rows swap
# positions (18768 <-> 31280) along with them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New "Periodo Mora" (column E) labels for rows 16-23.
$ws.Range("E16").Value = "1803"
$ws.Range("E17").Value = "2402"
$ws.Range("E18").Value = "2401"
$ws.Range("E19").Value = "2312"
$ws.Range("E20").Value = "2311"
$ws.Range("E21").Value = "2310"
$ws.Range("E22").Value = "2309"
$ws.Range("E23").Value = "2308"

# "Valor Mora" (column F) amounts follow the re-labelled periods.
$ws.Range("F16").Value = 31280
$ws.Range("F17").Value = 18768
$ws.Range("F18").Value = 31280
$ws.Range("F19").Value = 31280
$ws.Range("F20").Value = 31280
$ws.Range("F21").Value = 31280
$ws.Range("F22").Value = 31280
$ws.Range("F23").Value = 31280
